# Add a new "localdb" command-type column to the hidden '#system' sheet,
# wire up the new named range, and fix up all the named ranges that sit to
# the right of the newly inserted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a brand new column before N. Everything that used to live in
#    columns N..AC (macro, mail, number, pdf, rdbms, redis, sms, sound,
#    ssh, step, web, webalert, webcookie, ws, ws.async, xml) shifts one
#    column to the right, landing in O..AD.
# ---------------------------------------------------------------------
$ws.Columns("N").Insert()

# ---------------------------------------------------------------------
# 2) Insert a new row into column A only (the "target" list), just above
#    the current "macro" entry (row 14), so "localdb" can be slotted in
#    alphabetically between "json" and "macro". This only pushes column A
#    down - the rest of the columns are independent lists and must stay
#    put.
# ---------------------------------------------------------------------
$ws.Range("A14").Insert(-4121)

# ---------------------------------------------------------------------
# 3) Populate the new cells.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "localdb"

$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 4) Fix up the workbook-level named ranges. The Insert() calls above do
#    not shift defined names automatically, so they're repointed by hand.
# ---------------------------------------------------------------------
$wb.Names.Item("macro").RefersTo      = '=''#system''!$O$2:$O$4'
$wb.Names.Item("mail").RefersTo       = '=''#system''!$P$2:$P$2'
$wb.Names.Item("number").RefersTo     = '=''#system''!$Q$2:$Q$16'
$wb.Names.Item("pdf").RefersTo        = '=''#system''!$R$2:$R$16'
$wb.Names.Item("rdbms").RefersTo      = '=''#system''!$S$2:$S$7'
$wb.Names.Item("redis").RefersTo      = '=''#system''!$T$2:$T$10'
$wb.Names.Item("sms").RefersTo        = '=''#system''!$U$2:$U$2'
$wb.Names.Item("sound").RefersTo      = '=''#system''!$V$2:$V$5'
$wb.Names.Item("ssh").RefersTo        = '=''#system''!$W$2:$W$9'
$wb.Names.Item("step").RefersTo       = '=''#system''!$X$2:$X$4'
$wb.Names.Item("web").RefersTo        = '=''#system''!$Y$2:$Y$127'
$wb.Names.Item("webalert").RefersTo   = '=''#system''!$Z$2:$Z$8'
$wb.Names.Item("webcookie").RefersTo  = '=''#system''!$AA$2:$AA$8'
$wb.Names.Item("ws").RefersTo         = '=''#system''!$AB$2:$AB$17'
$wb.Names.Item("ws.async").RefersTo   = '=''#system''!$AC$2:$AC$8'
$wb.Names.Item("xml").RefersTo        = '=''#system''!$AD$2:$AD$21'

# "target" (column A) grew by one row because of the localdb insertion.
$wb.Names.Item("target").RefersTo     = '=''#system''!$A$2:$A$30'

# Brand new named range for the localdb command list.
$wb.Names.Add("localdb", '=''#system''!$N$2:$N$7')
